$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# --- Fix existing rows 48 and 49: CRM was actually opened 20210720, not 20210721 ---
$ws.Range("E48").Value = 180
$ws.Range("F48").Value = "CRM OPENED 20210720"

$ws.Range("E49").Value = 180
$ws.Range("F49").Value = "CRM OPENED 20210720"

# --- Add new data rows 50-52 for CBLS blue tank / quarantine tank TA DMBP runs ---

# Row 50
$ws.Range("A50").Value = 20211117
$ws.Range("B50").Value = 2231.686
$ws.Range("C50").Value = 2224.47
$ws.Range("D50").Formula = "=100*(B50-C50)/C50"
$ws.Range("E50").Value = 180
$ws.Range("F50").Value = "CRM OPENED 20210720"

# Row 51
$ws.Range("A51").Value = 20211206
$ws.Range("B51").Value = 2228.197
$ws.Range("C51").Value = 2224.47
$ws.Range("D51").Formula = "=100*(B51-C51)/C51"
$ws.Range("E51").Value = 180
$ws.Range("F51").Value = "CRM OPENED 20210720"

# Row 52 - new CRM opened 20220118
$ws.Range("A52").Value = 20220118
$ws.Range("B52").Value = 2203.673
$ws.Range("C52").Value = 2224.47
$ws.Range("D52").Formula = "=100*(B52-C52)/C52"
$ws.Range("E52").Value = 180
$ws.Range("F52").Value = "CRM OPENED 20220118"

# --- Update selection / active cell to reflect new last-entry position ---
$ws.Range("G52").Select()

$wb.Save()
